# Generate Report for Handoff
#
# Updates the localization-status workbook so the four "Ready for handoff"
# rows (6b8b06ff.., 8915f644.., a86c6719.., dd456874..) on both the zh-cn
# and de-de sheets reflect:
#   - Priority changed from "low" to "ht"
#   - Latest Handoff Datetime refreshed to the new handoff generation time

$wb = $excel.ActiveWorkbook

$rows = 4..7

$sheetUpdates = @{
    "zh-cn" = "2016-08-13 22:39:42"
    "de-de" = "2016-08-13 22:39:50"
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newHandoffDatetime = $sheetUpdates[$sheetName]

    foreach ($r in $rows) {
        $ws.Range("E$r").Value = "ht"
        $ws.Range("H$r").Value = $newHandoffDatetime
    }
}
